$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.86"
$ws.Range("E2").Value = "'1.28%"
$ws.Range("D3").Value = "'41.06"
$ws.Range("E3").Value = "'2.25%"
$ws.Range("D4").Value = "'5.746"
$ws.Range("E4").Value = "'-1.87%"
$ws.Range("D5").Value = "'0.08141"
$ws.Range("E5").Value = "'1.56%"
$ws.Range("D6").Value = "'2.046"
$ws.Range("E6").Value = "'6.76%"
$ws.Range("E7").Value = "'0.43%"
$ws.Range("D8").Value = "'4.508"
$ws.Range("E8").Value = "'-1.52%"
$ws.Range("E9").Value = "'0.79%"
$ws.Range("D10").Value = "'0.9207"
$ws.Range("E10").Value = "'-2.05%"
$ws.Range("D11").Value = "'0.1239"
$ws.Range("E11").Value = "'-0.73%"
$ws.Range("D12").Value = "'0.1947"
$ws.Range("E12").Value = "'-0.92%"
$ws.Range("D13").Value = "'8.295"
$ws.Range("E13").Value = "'-6.14%"
$ws.Range("D14").Value = "'0.09335"
$ws.Range("E14").Value = "'1.59%"
$ws.Range("D15").Value = "'0.03658"
$ws.Range("E15").Value = "'2.11%"
$ws.Range("E16").Value = "'9.55%"
$ws.Range("D17").Value = "'0.001303"
$ws.Range("E17").Value = "'-0.59%"
$ws.Range("D18").Value = "'0.006156"
$ws.Range("E18").Value = "'0.42%"
$ws.Range("D19").Value = "'3.382"
$ws.Range("E19").Value = "'0.80%"
$ws.Range("E20").Value = "'-1.20%"
$ws.Range("D21").Value = "'0.1417"
$ws.Range("E21").Value = "'-1.13%"
$ws.Range("E22").Value = "'9.57%"
$ws.Range("D23").Value = "'0.04436"
$ws.Range("E23").Value = "'0.04%"
$ws.Range("D24").Value = "'0.001259"
$ws.Range("E24").Value = "'-0.12%"
$ws.Range("D25").Value = "'0.004384"
$ws.Range("E25").Value = "'1.68%"
$ws.Range("E26").Value = "'8.49%"
$ws.Range("D39").Value = "'0.02785"
$ws.Range("E39").Value = "'15.46%"
$ws.Range("D40").Value = "'0.05494"
$ws.Range("E40").Value = "'4.43%"
$ws.Range("D41").Value = "'0.007598"
$ws.Range("E41").Value = "'1.49%"
$ws.Range("D42").Value = "'0.009938"
$ws.Range("D43").Value = "'0.1424"
$ws.Range("E43").Value = "'0.52%"
$ws.Range("D44").Value = "'0.002118"
$ws.Range("E44").Value = "'-0.43%"
$ws.Range("D45").Value = "'0.01189"
$ws.Range("E45").Value = "'11.95%"
$ws.Range("D46").Value = "'0.00006746"
$ws.Range("E46").Value = "'-1.19%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.43%"
$ws.Range("D48").Value = "'0.002939"
$ws.Range("E48").Value = "'-6.84%"
$ws.Range("D49").Value = "'0.002277"
$ws.Range("E49").Value = "'59.79%"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'-0.43%"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'-0.43%"
